$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.174.59"
$ws.Range("E2").Value = "  +2.10%  "

$ws.Range("D3").Value = "3.777.61"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "624.52"
$ws.Range("E5").Value = "  +4.25%  "

$ws.Range("D6").Value = "166.09"
$ws.Range("E6").Value = "  +1.81%  "

$ws.Range("D7").Value = "3.775.06"
$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("E10").Value = "  +2.76%  "

$ws.Range("D11").Value = "0.455"
$ws.Range("E11").Value = "  +2.40%  "

$ws.Range("D12").Value = "6.70"
$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("E13").Value = "  +0.91%  "

$ws.Range("D14").Value = "35.68"
$ws.Range("E14").Value = "  +1.61%  "

$ws.Range("D15").Value = "4.416.40"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("D16").Value = "3.783.06"
$ws.Range("E16").Value = "  +0.25%  "

$ws.Range("D17").Value = "69.262.06"
$ws.Range("E17").Value = "  +2.20%  "

$ws.Range("E18").Value = "  -2.60%  "

$ws.Range("E19").Value = "  +1.68%  "

$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("D21").Value = "467.96"
$ws.Range("E21").Value = "  +2.30%  "

$ws.Range("D22").Value = "9.61"
$ws.Range("E22").Value = "  +1.55%  "

$ws.Range("E23").Value = "  +2.13%  "

$ws.Range("D24").Value = "0.0000148"
$ws.Range("E24").Value = "  +3.97%  "

$ws.Range("D25").Value = "83.18"
$ws.Range("E25").Value = "  +0.59%  "

$ws.Range("D26").Value = "12.02"
$ws.Range("E26").Value = "  +1.60%  "

$ws.Range("E27").Value = "  +3.69%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.11%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "10.00"
$ws.Range("E29").Value = "  +1.71%  "

$ws.Range("D30").Value = "3.928.68"
$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("D31").Value = "2.67"
$ws.Range("E31").Value = "  +3.18%  "

$ws.Range("E32").Value = "  +2.35%  "

$ws.Range("D33").Value = "7.26"
$ws.Range("E33").Value = "  +1.09%  "

$ws.Range("D34").Value = "28.78"
$ws.Range("E34").Value = "  -0.26%  "

$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.166"
$ws.Range("E36").Value = "  +16.08%  "

$ws.Range("D37").Value = "3.730.29"
$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").Value = "9.00"
$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("E39").Value = "  +2.69%  "

$ws.Range("D40").Value = "3.41"
$ws.Range("E40").Value = "  +7.69%  "

$ws.Range("E41").Value = "  +0.76%  "

$ws.Range("D42").Value = "0.965"
$ws.Range("E42").Value = "  -1.15%  "

$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.298"
$ws.Range("E45").Value = "  +0.96%  "

$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").Value = "43.19"
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("E47").Value = "  +4.07%  "

$ws.Range("D48").Value = "46.69"
$ws.Range("E48").Value = "  -0.95%  "

$ws.Range("D49").Value = "151.71"
$ws.Range("E49").Value = "  -0.42%  "

$ws.Range("E50").Value = "  +1.89%  "

$ws.Range("E51").Value = "  +0.47%  "
